# Actualización automática 2025-06-03 09:25:07
# Populate the "PRESUPUESTO" (column G) budget figures on the
# "VENTA MENSUAL" sheet for the rows that received a budget assignment,
# and refresh the column total on the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$presupuesto = @{
    2  = 1500
    3  = 500
    4  = 1000
    8  = 1500
    9  = 500
    10 = 1500
    11 = 1000
    12 = 500
    13 = 1500
    15 = 500
    16 = 500
    17 = 1500
    20 = 500
    21 = 500
    24 = 300
    25 = 2000
    27 = 500
    28 = 500
}

foreach ($row in $presupuesto.Keys) {
    $ws.Cells.Item($row, 7).Value = $presupuesto[$row]
}

# Row 30 holds the column totals; update the PRESUPUESTO total (column G)
# to match the sum of the values just entered.
$ws.Cells.Item(30, 7).Value = 16300
